# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns AD, AE, AF on the header row.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header style (bold, centered, bordered) used by the
# rest of row 1 (e.g. AC1) so the new header cells look the same.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Every data row (2-44) gets the same constant team record values.
$ws.Range("AD2:AD44").Value = 94
$ws.Range("AE2:AE44").Value = 68
$ws.Range("AF2:AF44").Value = 0
